# Build site at 2022-09-26 16:07:08 UTC
# Simplify the LOM3108 syllabus sheet: the long free-text answers for
# Objetivos/Programa resumido/Programa/Metodo/Criterio/Bibliografia are
# removed from the sheet (site content trimmed during the build), which
# shifts every row below "Docentes responsaveis" up by one. Only the cells
# whose text actually changes are touched here, so untouched numeric-
# looking / date-looking text cells (credits, dates, etc.) keep their
# original text type instead of being re-interpreted by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C10").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C15").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas."
$ws.Range("C19").Value = "O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina."
$ws.Range("C20").Value = "A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "não há"
$ws.Range("C21").Value = "não há"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOM3104 -  Projeto Integrado em Engenharia de Materiais I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3104 -  Projeto Integrado em Engenharia de Materiais I  (Requisito fraco)`n"


# B19/B20 are brand-new cells (column B previously had no content on these
# rows). The sheet's <cols> entries for column B overlap (min=1-2 style=1,
# then min=2 style=2), so a freshly created B-cell picks up the first
# (wrong, bold) style instead of the intended wrapped-text style used by
# every other column-B cell. Copy the correct format from B10 to fix it.
$ws.Range("B10").Copy()
$ws.Range("B19:B20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Clear cells that no longer hold content in the updated layout
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()


# Row heights for the affected rows (13-23)
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30


# Remove now-unused trailing rows 24-26 (content moved up into rows 13-23)
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()
